$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 with the new venue details (Melbourne Central Lion Hotel)
$ws.Range("A2").Value = "Melbourne"
$ws.Range("B2").Value = "Melbourne Central Lion Hotel, 211 La Trobe Street"
$ws.Range("C2").Value = "28/12/2020 10:30pm-12.00am"
$ws.Range("D2").Value = "Case attended venue"
$ws.Range("E2").Value = "new"

# Add row 3 - superseded exposure window for the same venue
$ws.Range("A3").Value = "Melbourne"
$ws.Range("B3").Value = "Melbourne Central Lion Hotel, 211 La Trobe Street"
$ws.Range("C3").Value = "28/12/2020 10pm - 12.30am"
$ws.Range("D3").Value = "Case attended venue"
$ws.Range("E3").Value = "old"

# Add row 4 - new Nando's exposure site
$ws.Range("A4").Value = "Melbourne"
$ws.Range("B4").Value = "Nandos  27 Elizabeth Street, Melbourne"
$ws.Range("C4").Value = "01/01/2021 1:00am - 2:00am"
$ws.Range("D4").Value = "Case dined at venue"
$ws.Range("E4").Value = "new"

# Autofit columns to match Excel's recalculated best-fit widths
$ws.Columns.AutoFit() | Out-Null
$ws.Range("A1").EntireColumn.ColumnWidth = 8.5
$ws.Range("B1").EntireColumn.ColumnWidth = 40
$ws.Range("C1").EntireColumn.ColumnWidth = 24.6666666666667
$ws.Range("D1").EntireColumn.ColumnWidth = 16.5

# Restore active selection to B3 as the last edited cell
$ws.Range("B3").Select() | Out-Null
